# The edit rotates the data of rows 2, 3 and 4 (all columns A:AY) upward
# by one position in a cycle: row 2 <- row 3, row 3 <- row 4, row 4 <- row 2
# (i.e. new_row2 = old_row4's data, new_row3 = old_row2's data,
#  new_row4 = old_row3's data). Record ids in column A / B stay attached
# to their own row's other data, so the whole row moves together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns used in the data rows (A through AY, matching the header row).
$cols = @(
  "A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z",
  "AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY"
)

# Snapshot the current ("before") values of rows 2-4 for every column first,
# so writes to one row don't affect values still to be read from another.
$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("${col}2").Value2
    $row3[$col] = $ws.Range("${col}3").Value2
    $row4[$col] = $ws.Range("${col}4").Value2
}

# Apply the cyclic rotation: new2 = old4, new3 = old2, new4 = old3.
foreach ($col in $cols) {
    $ws.Range("${col}2").Value = $row4[$col]
    $ws.Range("${col}3").Value = $row2[$col]
    $ws.Range("${col}4").Value = $row3[$col]
}
